$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 348.46155
$ws.Range("I18").Value = 348.46155
$ws.Range("K18").Value = 348.46155
$ws.Range("M18").Value = -64.46154999999999
$ws.Range("H33").Value = 270.45834
$ws.Range("I33").Value = 216.14285
$ws.Range("J33").Value = 650.6667
$ws.Range("K33").Value = 216.14285
$ws.Range("L33").Value = 650.6667
$ws.Range("M33").Value = 12.85714999999999
$ws.Range("N33").Value = -1108.6667
$ws.Range("H64").Value = 411986.3
$ws.Range("J64").Value = 4366.5
$ws.Range("L64").Value = 4366.5
$ws.Range("N64").Value = -4862.5
$ws.Range("H67").Value = 411986.3
$ws.Range("J67").Value = 4366.5
$ws.Range("L67").Value = 4366.5
$ws.Range("N67").Value = -6082.5
$ws.Range("H112").Value = 1937.7
$ws.Range("J112").Value = 2425.2856
$ws.Range("L112").Value = 7275.8568
$ws.Range("N112").Value = -9491.856800000001
$ws.Range("H113").Value = 2521.389
$ws.Range("I113").Value = 2012.3572
$ws.Range("J113").Value = 4303
$ws.Range("K113").Value = 2012.3572
$ws.Range("L113").Value = 4303
$ws.Range("M113").Value = 1241.6428
$ws.Range("N113").Value = -10811
$ws.Range("H129").Value = 1083.2632
$ws.Range("J129").Value = 1126.7778
$ws.Range("L129").Value = 3380.3334
$ws.Range("N129").Value = -13380.3334
$ws.Range("H138").Value = 5624.528
$ws.Range("I138").Value = 1854.6666
$ws.Range("J138").Value = 7176.8237
$ws.Range("K138").Value = 5563.9998
$ws.Range("L138").Value = 21530.4711
$ws.Range("M138").Value = -423.9997999999996
$ws.Range("N138").Value = -31810.4711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1891.3182
$ws.Range("I2").Value = 1942.6316
$ws.Range("J2").Value = 1566.3334
$ws.Range("K2").Value = 1942.6316
$ws.Range("L2").Value = 1566.3334
$ws.Range("M2").Value = -1829.6316
$ws.Range("N2").Value = -1792.3334
$ws.Range("H45").Value = 1464.579
$ws.Range("I45").Value = 1467.1875
$ws.Range("J45").Value = 1450.6666
$ws.Range("K45").Value = 1467.1875
$ws.Range("L45").Value = 1450.6666
$ws.Range("M45").Value = -1090.1875
$ws.Range("N45").Value = -2204.6666
$ws.Range("H61").Value = 5584.9375
$ws.Range("I61").Value = 3624.2778
$ws.Range("J61").Value = 11466.917
$ws.Range("K61").Value = 3624.2778
$ws.Range("L61").Value = 11466.917
$ws.Range("M61").Value = -3412.2778
$ws.Range("N61").Value = -11890.917
$ws.Range("H68").Value = 28929.334
$ws.Range("I68").Value = 10590
$ws.Range("K68").Value = 10590
$ws.Range("M68").Value = -9779
$ws.Range("H71").Value = 28929.334
$ws.Range("I71").Value = 10590
$ws.Range("K71").Value = 31770
$ws.Range("M71").Value = -27714
$ws.Range("H74").Value = 8521.267
$ws.Range("I74").Value = 6978.4
$ws.Range("J74").Value = 11607
$ws.Range("K74").Value = 6978.4
$ws.Range("L74").Value = 11607
$ws.Range("M74").Value = -6104.4
$ws.Range("N74").Value = -13355
$ws.Range("H77").Value = 8521.267
$ws.Range("I77").Value = 6978.4
$ws.Range("J77").Value = 11607
$ws.Range("K77").Value = 34892
$ws.Range("L77").Value = 58035
$ws.Range("M77").Value = -30524
$ws.Range("N77").Value = -66771
$ws.Range("H110").Value = 1580.0625
$ws.Range("I110").Value = 1561.4615
$ws.Range("J110").Value = 1660.6666
$ws.Range("K110").Value = 1561.4615
$ws.Range("L110").Value = 1660.6666
$ws.Range("M110").Value = 483.5385000000001
$ws.Range("N110").Value = -5750.6666
$ws.Range("H116").Value = 1891.3182
$ws.Range("I116").Value = 1942.6316
$ws.Range("J116").Value = 1566.3334
$ws.Range("K116").Value = 1942.6316
$ws.Range("L116").Value = 1566.3334
$ws.Range("M116").Value = 351.3684000000001
$ws.Range("N116").Value = -6154.3334
$ws.Range("H132").Value = 1861.381
$ws.Range("I132").Value = 1183.4839
$ws.Range("J132").Value = 3771.818
$ws.Range("K132").Value = 3550.4517
$ws.Range("L132").Value = 11315.454
$ws.Range("M132").Value = -1020.4517
$ws.Range("N132").Value = -16375.454
$ws.Range("H136").Value = 5584.9375
$ws.Range("I136").Value = 3624.2778
$ws.Range("J136").Value = 11466.917
$ws.Range("K136").Value = 10872.8334
$ws.Range("L136").Value = 34400.751
$ws.Range("M136").Value = -8322.8334
$ws.Range("N136").Value = -39500.751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1891.3182
$ws.Range("I3").Value = 1942.6316
$ws.Range("J3").Value = 1566.3334
$ws.Range("K3").Value = 1942.6316
$ws.Range("L3").Value = 1566.3334
$ws.Range("M3").Value = -1828.6316
$ws.Range("N3").Value = -1794.3334
$ws.Range("H36").Value = 3356.1667
$ws.Range("I36").Value = 884.25
$ws.Range("K36").Value = 884.25
$ws.Range("M36").Value = -350.25
$ws.Range("H105").Value = 5389.45
$ws.Range("I105").Value = 4940.75
$ws.Range("J105").Value = 6062.5
$ws.Range("K105").Value = 4940.75
$ws.Range("L105").Value = 6062.5
$ws.Range("M105").Value = -3193.75
$ws.Range("N105").Value = -9556.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4290.7173
$ws.Range("I31").Value = 4893.4443
$ws.Range("J31").Value = 3434.2104
$ws.Range("K31").Value = 4893.4443
$ws.Range("L31").Value = 3434.2104
$ws.Range("M31").Value = -4598.4443
$ws.Range("N31").Value = -4024.2104
$ws.Range("H34").Value = 4290.7173
$ws.Range("I34").Value = 4893.4443
$ws.Range("J34").Value = 3434.2104
$ws.Range("K34").Value = 4893.4443
$ws.Range("L34").Value = 3434.2104
$ws.Range("M34").Value = -4691.4443
$ws.Range("N34").Value = -3838.2104
$ws.Range("H86").Value = 2264.72
$ws.Range("I86").Value = 2700.4375
$ws.Range("K86").Value = 2700.4375
$ws.Range("M86").Value = -1577.4375
$ws.Range("H89").Value = 2264.72
$ws.Range("I89").Value = 2700.4375
$ws.Range("K89").Value = 13502.1875
$ws.Range("M89").Value = -7886.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 52.75
$ws.Range("I2").Value = 23.333334
$ws.Range("K2").Value = 140.000004
$ws.Range("M2").Value = -27.00000399999999
$ws.Range("H123").Value = 2213.7144
$ws.Range("I123").Value = 350
$ws.Range("J123").Value = 2959.2
$ws.Range("K123").Value = 1050
$ws.Range("L123").Value = 8877.599999999999
$ws.Range("M123").Value = 1400
$ws.Range("N123").Value = -13777.6
$ws.Range("H131").Value = 31191.094
$ws.Range("I131").Value = 1276.9524
$ws.Range("J131").Value = 88299.91
$ws.Range("K131").Value = 3830.857199999999
$ws.Range("L131").Value = 264899.73
$ws.Range("M131").Value = 1209.142800000001
$ws.Range("N131").Value = -274979.73
$ws.Range("H136").Value = 2836.6072
$ws.Range("I136").Value = 1181.6666
$ws.Range("J136").Value = 3287.9546
$ws.Range("K136").Value = 3544.9998
$ws.Range("L136").Value = 9863.863799999999
$ws.Range("M136").Value = 1555.0002
$ws.Range("N136").Value = -20063.8638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H70").Value = 5303.6274
$ws.Range("I70").Value = 4736.875
$ws.Range("K70").Value = 4736.875
$ws.Range("M70").Value = -4466.875
$ws.Range("H73").Value = 5303.6274
$ws.Range("I73").Value = 4736.875
$ws.Range("K73").Value = 4736.875
$ws.Range("M73").Value = -3800.875
$ws.Range("H80").Value = 10761.23
$ws.Range("I80").Value = 20725
$ws.Range("J80").Value = 6332.8887
$ws.Range("K80").Value = 20725
$ws.Range("L80").Value = 6332.8887
$ws.Range("M80").Value = -19727
$ws.Range("N80").Value = -8328.8887
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 10761.23
$ws.Range("I83").Value = 20725
$ws.Range("J83").Value = 6332.8887
$ws.Range("K83").Value = 103625
$ws.Range("L83").Value = 31664.4435
$ws.Range("M83").Value = -98633
$ws.Range("N83").Value = -41648.4435
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 1793.8695
$ws.Range("I132").Value = 1685.7222
$ws.Range("J132").Value = 2183.2
$ws.Range("K132").Value = 5057.1666
$ws.Range("L132").Value = 6549.599999999999
$ws.Range("M132").Value = -2527.1666
$ws.Range("N132").Value = -11609.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5799.933
$ws.Range("I7").Value = 5167
$ws.Range("J7").Value = 8331.666999999999
$ws.Range("K7").Value = 5167
$ws.Range("L7").Value = 8331.666999999999
$ws.Range("M7").Value = -5055
$ws.Range("N7").Value = -8555.666999999999
$ws.Range("H93").Value = 500
$ws.Range("I93").Value = 500
$ws.Range("K93").Value = 500
$ws.Range("M93").Value = 748
$ws.Range("H126").Value = 5799.933
$ws.Range("I126").Value = 5167
$ws.Range("J126").Value = 8331.666999999999
$ws.Range("K126").Value = 15501
$ws.Range("L126").Value = 24995.001
$ws.Range("M126").Value = -13031
$ws.Range("N126").Value = -29935.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 27842.6
$ws.Range("J75").Value = 31053.25
$ws.Range("L75").Value = 31053.25
$ws.Range("N75").Value = -32925.25
$ws.Range("H78").Value = 27842.6
$ws.Range("J78").Value = 31053.25
$ws.Range("L78").Value = 93159.75
$ws.Range("N78").Value = -102519.75
$ws.Range("H86").Value = 49000
$ws.Range("J86").Value = 49000
$ws.Range("L86").Value = 49000
$ws.Range("N86").Value = -51246
$ws.Range("H89").Value = 49000
$ws.Range("J89").Value = 49000
$ws.Range("L89").Value = 245000
$ws.Range("N89").Value = -256232
$ws.Range("H126").Value = 1449.75
$ws.Range("I126").Value = 1493.125
$ws.Range("J126").Value = 1276.25
$ws.Range("K126").Value = 4479.375
$ws.Range("L126").Value = 3828.75
$ws.Range("M126").Value = -2009.375
$ws.Range("N126").Value = -8768.75
$ws.Range("H132").Value = 1841.7843
$ws.Range("I132").Value = 847.0857
$ws.Range("J132").Value = 4017.6875
$ws.Range("K132").Value = 2541.2571
$ws.Range("L132").Value = 12053.0625
$ws.Range("M132").Value = -11.25709999999981
$ws.Range("N132").Value = -17113.0625

Write-Host "Edit complete"